$wb = $excel.ActiveWorkbook

# Update Proveedor sheet values
$ws = $wb.Worksheets.Item("Proveedor")
$ws.Range("B2").Value = "Proveedor nuevo nuevo"
$ws.Range("B3").Value = "20.630.735-8"
$ws.Range("B4").Value = "t2est@proveedorejemplo.com"
$ws.Range("B5").Value = "'387654321"
$ws.Range("B6").Value = "Called"

# Remove Productos sheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Productos").Delete()
